$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 51.22717533333333
$ws.Cells.Item(2, 8).Value = 153.681526
$ws.Cells.Item(2, 9).Value = 0.1770805608477903
$ws.Cells.Item(2, 10).Value = 0.1770805608477904
$ws.Cells.Item(2, 13).Value = 21.85308466666666
$ws.Cells.Item(2, 14).Value = 65.559254
$ws.Cells.Item(2, 15).Value = 0.407053040353553
$ws.Cells.Item(2, 16).Value = 0.407053040353553
$ws.Cells.Item(2, 17).Value = 1119.471799793511
$ws.Cells.Item(2, 18).Value = 10075.2461981416
$ws.Cells.Item(2, 19).Value = 0.0720811806806054
$ws.Cells.Item(2, 20).Value = 0.07208118068060541
$ws.Cells.Item(3, 7).Value = 51.22717533333333
$ws.Cells.Item(3, 8).Value = 153.681526
$ws.Cells.Item(3, 9).Value = 0.1770805608477903
$ws.Cells.Item(3, 10).Value = 0.1770805608477904
$ws.Cells.Item(3, 15).Value = 0.1342711086924142
$ws.Cells.Item(3, 16).Value = 0.1342711086924142
$ws.Cells.Item(3, 17).Value = 369.2705981943051
$ws.Cells.Item(3, 18).Value = 3323.435383748746
$ws.Cells.Item(3, 19).Value = 0.02377680323290732
$ws.Cells.Item(3, 20).Value = 0.02377680323290733
$ws.Cells.Item(4, 7).Value = 51.22717533333333
$ws.Cells.Item(4, 8).Value = 153.681526
$ws.Cells.Item(4, 9).Value = 0.1770805608477903
$ws.Cells.Item(4, 10).Value = 0.1770805608477904
$ws.Cells.Item(4, 13).Value = 11.375406
$ws.Cells.Item(4, 14).Value = 34.126218
$ws.Cells.Item(4, 15).Value = 0.2118874139822907
$ws.Cells.Item(4, 16).Value = 0.2118874139822907
$ws.Cells.Item(4, 17).Value = 582.7299176498519
$ws.Cells.Item(4, 18).Value = 5244.569258848668
$ws.Cells.Item(4, 19).Value = 0.03752114210457197
$ws.Cells.Item(4, 20).Value = 0.03752114210457198
$ws.Cells.Item(5, 7).Value = 51.22717533333333
$ws.Cells.Item(5, 8).Value = 153.681526
$ws.Cells.Item(5, 9).Value = 0.1770805608477903
$ws.Cells.Item(5, 10).Value = 0.1770805608477904
$ws.Cells.Item(5, 13).Value = 3.401340666666667
$ws.Cells.Item(5, 14).Value = 10.204022
$ws.Cells.Item(5, 15).Value = 0.06335609277882483
$ws.Cells.Item(5, 16).Value = 0.06335609277882483
$ws.Cells.Item(5, 17).Value = 174.2410746997302
$ws.Cells.Item(5, 18).Value = 1568.169672297572
$ws.Cells.Item(5, 19).Value = 0.01121913244239894
$ws.Cells.Item(5, 20).Value = 0.01121913244239894
$ws.Cells.Item(6, 7).Value = 51.22717533333333
$ws.Cells.Item(6, 8).Value = 153.681526
$ws.Cells.Item(6, 9).Value = 0.1770805608477903
$ws.Cells.Item(6, 10).Value = 0.1770805608477904
$ws.Cells.Item(6, 13).Value = 9.847764666666666
$ws.Cells.Item(6, 14).Value = 29.543294
$ws.Cells.Item(6, 15).Value = 0.1834323441929172
$ws.Cells.Item(6, 16).Value = 0.1834323441929172
$ws.Cells.Item(6, 17).Value = 504.4731672207382
$ws.Cells.Item(6, 18).Value = 4540.258504986644
$ws.Cells.Item(6, 19).Value = 0.03248230238730668
$ws.Cells.Item(6, 20).Value = 0.03248230238730669
$ws.Cells.Item(7, 9).Value = 0.2662631264141754
$ws.Cells.Item(7, 10).Value = 0.2662631264141754
$ws.Cells.Item(7, 13).Value = 21.85308466666666
$ws.Cells.Item(7, 14).Value = 65.559254
$ws.Cells.Item(7, 15).Value = 0.407053040353553
$ws.Cells.Item(7, 16).Value = 0.407053040353553
$ws.Cells.Item(7, 17).Value = 1683.268112086757
$ws.Cells.Item(7, 18).Value = 15149.41300878082
$ws.Cells.Item(7, 19).Value = 0.1083832151409325
$ws.Cells.Item(7, 20).Value = 0.1083832151409325
$ws.Cells.Item(8, 9).Value = 0.2662631264141754
$ws.Cells.Item(8, 10).Value = 0.2662631264141754
$ws.Cells.Item(8, 15).Value = 0.1342711086924142
$ws.Cells.Item(8, 16).Value = 0.1342711086924142
$ws.Cells.Item(8, 19).Value = 0.03575144518753977
$ws.Cells.Item(8, 20).Value = 0.03575144518753978
$ws.Cells.Item(9, 9).Value = 0.2662631264141754
$ws.Cells.Item(9, 10).Value = 0.2662631264141754
$ws.Cells.Item(9, 13).Value = 11.375406
$ws.Cells.Item(9, 14).Value = 34.126218
$ws.Cells.Item(9, 15).Value = 0.2118874139822907
$ws.Cells.Item(9, 16).Value = 0.2118874139822907
$ws.Cells.Item(9, 17).Value = 876.208483786608
$ws.Cells.Item(9, 18).Value = 7885.876354079472
$ws.Cells.Item(9, 19).Value = 0.05641780529473939
$ws.Cells.Item(9, 20).Value = 0.05641780529473939
$ws.Cells.Item(10, 9).Value = 0.2662631264141754
$ws.Cells.Item(10, 10).Value = 0.2662631264141754
$ws.Cells.Item(10, 13).Value = 3.401340666666667
$ws.Cells.Item(10, 14).Value = 10.204022
$ws.Cells.Item(10, 15).Value = 0.06335609277882483
$ws.Cells.Item(10, 16).Value = 0.06335609277882483
$ws.Cells.Item(10, 17).Value = 261.9935981521653
$ws.Cells.Item(10, 18).Value = 2357.942383369488
$ws.Cells.Item(10, 19).Value = 0.01686939134067646
$ws.Cells.Item(10, 20).Value = 0.01686939134067646
$ws.Cells.Item(11, 9).Value = 0.2662631264141754
$ws.Cells.Item(11, 10).Value = 0.2662631264141754
$ws.Cells.Item(11, 13).Value = 9.847764666666666
$ws.Cells.Item(11, 14).Value = 29.543294
$ws.Cells.Item(11, 15).Value = 0.1834323441929172
$ws.Cells.Item(11, 16).Value = 0.1834323441929172
$ws.Cells.Item(11, 17).Value = 758.5395147449973
$ws.Cells.Item(11, 18).Value = 6826.855632704976
$ws.Cells.Item(11, 19).Value = 0.04884126945028723
$ws.Cells.Item(11, 20).Value = 0.04884126945028724
$ws.Cells.Item(12, 7).Value = 72.76991766666667
$ws.Cells.Item(12, 8).Value = 218.309753
$ws.Cells.Item(12, 9).Value = 0.2515488654100336
$ws.Cells.Item(12, 10).Value = 0.2515488654100336
$ws.Cells.Item(12, 13).Value = 21.85308466666666
$ws.Cells.Item(12, 14).Value = 65.559254
$ws.Cells.Item(12, 15).Value = 0.407053040353553
$ws.Cells.Item(12, 16).Value = 0.407053040353553
$ws.Cells.Item(12, 17).Value = 1590.247171956029
$ws.Cells.Item(12, 18).Value = 14312.22454760426
$ws.Cells.Item(12, 19).Value = 0.1023937304626409
$ws.Cells.Item(12, 20).Value = 0.1023937304626409
$ws.Cells.Item(13, 7).Value = 72.76991766666667
$ws.Cells.Item(13, 8).Value = 218.309753
$ws.Cells.Item(13, 9).Value = 0.2515488654100336
$ws.Cells.Item(13, 10).Value = 0.2515488654100336
$ws.Cells.Item(13, 15).Value = 0.1342711086924142
$ws.Cells.Item(13, 16).Value = 0.1342711086924142
$ws.Cells.Item(13, 17).Value = 524.5612480576293
$ws.Cells.Item(13, 18).Value = 4721.051232518663
$ws.Cells.Item(13, 19).Value = 0.0337757450489241
$ws.Cells.Item(13, 20).Value = 0.0337757450489241
$ws.Cells.Item(14, 7).Value = 72.76991766666667
$ws.Cells.Item(14, 8).Value = 218.309753
$ws.Cells.Item(14, 9).Value = 0.2515488654100336
$ws.Cells.Item(14, 10).Value = 0.2515488654100336
$ws.Cells.Item(14, 13).Value = 11.375406
$ws.Cells.Item(14, 14).Value = 34.126218
$ws.Cells.Item(14, 15).Value = 0.2118874139822907
$ws.Cells.Item(14, 16).Value = 0.2118874139822907
$ws.Cells.Item(14, 17).Value = 827.787358044906
$ws.Cells.Item(14, 18).Value = 7450.086222404155
$ws.Cells.Item(14, 19).Value = 0.05330003858191133
$ws.Cells.Item(14, 20).Value = 0.05330003858191132
$ws.Cells.Item(15, 7).Value = 72.76991766666667
$ws.Cells.Item(15, 8).Value = 218.309753
$ws.Cells.Item(15, 9).Value = 0.2515488654100336
$ws.Cells.Item(15, 10).Value = 0.2515488654100336
$ws.Cells.Item(15, 13).Value = 3.401340666666667
$ws.Cells.Item(15, 14).Value = 10.204022
$ws.Cells.Item(15, 15).Value = 0.06335609277882483
$ws.Cells.Item(15, 16).Value = 0.06335609277882483
$ws.Cells.Item(15, 17).Value = 247.5152802696184
$ws.Cells.Item(15, 18).Value = 2227.637522426566
$ws.Cells.Item(15, 19).Value = 0.01593715325532621
$ws.Cells.Item(15, 20).Value = 0.01593715325532621
$ws.Cells.Item(16, 7).Value = 72.76991766666667
$ws.Cells.Item(16, 8).Value = 218.309753
$ws.Cells.Item(16, 9).Value = 0.2515488654100336
$ws.Cells.Item(16, 10).Value = 0.2515488654100336
$ws.Cells.Item(16, 13).Value = 9.847764666666666
$ws.Cells.Item(16, 14).Value = 29.543294
$ws.Cells.Item(16, 15).Value = 0.1834323441929172
$ws.Cells.Item(16, 16).Value = 0.1834323441929172
$ws.Cells.Item(16, 17).Value = 716.6210239940425
$ws.Cells.Item(16, 18).Value = 6449.589215946382
$ws.Cells.Item(16, 19).Value = 0.04614219806123108
$ws.Cells.Item(16, 20).Value = 0.04614219806123108
$ws.Cells.Item(17, 7).Value = 32.02005133333333
$ws.Cells.Item(17, 8).Value = 96.06015400000001
$ws.Cells.Item(17, 9).Value = 0.110685951579145
$ws.Cells.Item(17, 10).Value = 0.110685951579145
$ws.Cells.Item(17, 13).Value = 21.85308466666666
$ws.Cells.Item(17, 14).Value = 65.559254
$ws.Cells.Item(17, 15).Value = 0.407053040353553
$ws.Cells.Item(17, 16).Value = 0.407053040353553
$ws.Cells.Item(17, 17).Value = 699.7368928183462
$ws.Cells.Item(17, 18).Value = 6297.632035365116
$ws.Cells.Item(17, 19).Value = 0.04505505311471712
$ws.Cells.Item(17, 20).Value = 0.04505505311471713
$ws.Cells.Item(18, 7).Value = 32.02005133333333
$ws.Cells.Item(18, 8).Value = 96.06015400000001
$ws.Cells.Item(18, 9).Value = 0.110685951579145
$ws.Cells.Item(18, 10).Value = 0.110685951579145
$ws.Cells.Item(18, 15).Value = 0.1342711086924142
$ws.Cells.Item(18, 16).Value = 0.1342711086924142
$ws.Cells.Item(18, 17).Value = 230.8162305091705
$ws.Cells.Item(18, 18).Value = 2077.346074582534
$ws.Cells.Item(18, 19).Value = 0.01486192543520667
$ws.Cells.Item(18, 20).Value = 0.01486192543520668
$ws.Cells.Item(19, 7).Value = 32.02005133333333
$ws.Cells.Item(19, 8).Value = 96.06015400000001
$ws.Cells.Item(19, 9).Value = 0.110685951579145
$ws.Cells.Item(19, 10).Value = 0.110685951579145
$ws.Cells.Item(19, 13).Value = 11.375406
$ws.Cells.Item(19, 14).Value = 34.126218
$ws.Cells.Item(19, 15).Value = 0.2118874139822907
$ws.Cells.Item(19, 16).Value = 0.2118874139822907
$ws.Cells.Item(19, 17).Value = 364.241084057508
$ws.Cells.Item(19, 18).Value = 3278.169756517573
$ws.Cells.Item(19, 19).Value = 0.02345296004427408
$ws.Cells.Item(19, 20).Value = 0.02345296004427408
$ws.Cells.Item(20, 7).Value = 32.02005133333333
$ws.Cells.Item(20, 8).Value = 96.06015400000001
$ws.Cells.Item(20, 9).Value = 0.110685951579145
$ws.Cells.Item(20, 10).Value = 0.110685951579145
$ws.Cells.Item(20, 13).Value = 3.401340666666667
$ws.Cells.Item(20, 14).Value = 10.204022
$ws.Cells.Item(20, 15).Value = 0.06335609277882483
$ws.Cells.Item(20, 16).Value = 0.06335609277882483
$ws.Cells.Item(20, 17).Value = 108.9111027488209
$ws.Cells.Item(20, 18).Value = 980.1999247393882
$ws.Cells.Item(20, 19).Value = 0.007012629417560823
$ws.Cells.Item(20, 20).Value = 0.007012629417560824
$ws.Cells.Item(21, 7).Value = 32.02005133333333
$ws.Cells.Item(21, 8).Value = 96.06015400000001
$ws.Cells.Item(21, 9).Value = 0.110685951579145
$ws.Cells.Item(21, 10).Value = 0.110685951579145
$ws.Cells.Item(21, 13).Value = 9.847764666666666
$ws.Cells.Item(21, 14).Value = 29.543294
$ws.Cells.Item(21, 15).Value = 0.1834323441929172
$ws.Cells.Item(21, 16).Value = 0.1834323441929172
$ws.Cells.Item(21, 17).Value = 315.3259301452529
$ws.Cells.Item(21, 18).Value = 2837.933371307277
$ws.Cells.Item(21, 19).Value = 0.02030338356738629
$ws.Cells.Item(21, 20).Value = 0.02030338356738629
$ws.Cells.Item(22, 7).Value = 56.24368933333333
$ws.Cells.Item(22, 8).Value = 168.731068
$ws.Cells.Item(22, 9).Value = 0.1944214957488556
$ws.Cells.Item(22, 10).Value = 0.1944214957488557
$ws.Cells.Item(22, 13).Value = 21.85308466666666
$ws.Cells.Item(22, 14).Value = 65.559254
$ws.Cells.Item(22, 15).Value = 0.407053040353553
$ws.Cells.Item(22, 16).Value = 0.407053040353553
$ws.Cells.Item(22, 17).Value = 1229.09810496703
$ws.Cells.Item(22, 18).Value = 11061.88294470327
$ws.Cells.Item(22, 19).Value = 0.07913986095465707
$ws.Cells.Item(22, 20).Value = 0.07913986095465708
$ws.Cells.Item(23, 7).Value = 56.24368933333333
$ws.Cells.Item(23, 8).Value = 168.731068
$ws.Cells.Item(23, 9).Value = 0.1944214957488556
$ws.Cells.Item(23, 10).Value = 0.1944214957488557
$ws.Cells.Item(23, 15).Value = 0.1342711086924142
$ws.Cells.Item(23, 16).Value = 0.1342711086924142
$ws.Cells.Item(23, 17).Value = 405.4320908703364
$ws.Cells.Item(23, 18).Value = 3648.888817833028
$ws.Cells.Item(23, 19).Value = 0.02610518978783634
$ws.Cells.Item(23, 20).Value = 0.02610518978783635
$ws.Cells.Item(24, 7).Value = 56.24368933333333
$ws.Cells.Item(24, 8).Value = 168.731068
$ws.Cells.Item(24, 9).Value = 0.1944214957488556
$ws.Cells.Item(24, 10).Value = 0.1944214957488557
$ws.Cells.Item(24, 13).Value = 11.375406
$ws.Cells.Item(24, 14).Value = 34.126218
$ws.Cells.Item(24, 15).Value = 0.2118874139822907
$ws.Cells.Item(24, 16).Value = 0.2118874139822907
$ws.Cells.Item(24, 17).Value = 639.794801104536
$ws.Cells.Item(24, 18).Value = 5758.153209940824
$ws.Cells.Item(24, 19).Value = 0.04119546795679395
$ws.Cells.Item(24, 20).Value = 0.04119546795679395
$ws.Cells.Item(25, 7).Value = 56.24368933333333
$ws.Cells.Item(25, 8).Value = 168.731068
$ws.Cells.Item(25, 9).Value = 0.1944214957488556
$ws.Cells.Item(25, 10).Value = 0.1944214957488557
$ws.Cells.Item(25, 13).Value = 3.401340666666667
$ws.Cells.Item(25, 14).Value = 10.204022
$ws.Cells.Item(25, 15).Value = 0.06335609277882483
$ws.Cells.Item(25, 16).Value = 0.06335609277882483
$ws.Cells.Item(25, 17).Value = 191.3039477728329
$ws.Cells.Item(25, 18).Value = 1721.735529955496
$ws.Cells.Item(25, 19).Value = 0.0123177863228624
$ws.Cells.Item(25, 20).Value = 0.0123177863228624
$ws.Cells.Item(26, 7).Value = 56.24368933333333
$ws.Cells.Item(26, 8).Value = 168.731068
$ws.Cells.Item(26, 9).Value = 0.1944214957488556
$ws.Cells.Item(26, 10).Value = 0.1944214957488557
$ws.Cells.Item(26, 13).Value = 9.847764666666666
$ws.Cells.Item(26, 14).Value = 29.543294
$ws.Cells.Item(26, 15).Value = 0.1834323441929172
$ws.Cells.Item(26, 16).Value = 0.1834323441929172
$ws.Cells.Item(26, 17).Value = 553.8746165397769
$ws.Cells.Item(26, 18).Value = 4984.871548857992
$ws.Cells.Item(26, 19).Value = 0.03566319072670587
$ws.Cells.Item(26, 20).Value = 0.03566319072670587
